$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the missing "people_score" value into column K for the rows
# --- that were shifted left by one column (K held "overview" text, L held
# --- "consensus" text, M held the trailing flag). Shift L,M -> M,N and
# --- drop the new numeric score into K.
function Insert-PeopleScore([int]$row, [double]$score) {
    $colM = $ws.Range("M" + $row).Value2
    $colL = $ws.Range("L" + $row).Value2
    $colK = $ws.Range("K" + $row).Value2
    $ws.Range("N" + $row).Value = $colM
    $ws.Range("M" + $row).Value = $colL
    $ws.Range("L" + $row).Value = $colK
    $ws.Range("K" + $row).Value = $score
}

Insert-PeopleScore 8 88
Insert-PeopleScore 11 76
Insert-PeopleScore 22 46
Insert-PeopleScore 25 90
Insert-PeopleScore 26 53

# --- Column B width ---
$ws.Columns("B").ColumnWidth = 23 + 2/3

# --- Sheet view: zoom, scroll position, selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 181
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K26").Select()
